$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Change Log")

# Update A13 with the date serial value 45942 (2025-10-12), matching style/format
# already applied to the cell (date number format).
$ws.Range("A13").Value = 45942

# Update B13 text to reflect completed testing of all instructions.
$newText = "Changes`n- MODIFIED: MyMIF.mif`n- COMPLETED: Testing (simulation) for all instructions currently in the file. Added comments about said testing.                                                                                                                                                                                                                                       "
$ws.Range("B13").Value = $newText

# Update the sheet view: scroll position and active selection.
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("A14").Select()
